# Update cryptocurrency Price (D) and Volume(1h) (E) columns for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.496.51"
$ws.Range("E2").Value = "  +0.78%  "

$ws.Range("D3").Value = "'2.014.88"
$ws.Range("E3").Value = "  -0.93%  "

$ws.Range("E4").Value = "  -0.25%  "

$ws.Range("D5").Value = "'252.19"
$ws.Range("E5").Value = "  +3.29%  "

$ws.Range("D6").Value = "'0.639"
$ws.Range("E6").Value = "  -2.83%  "

$ws.Range("D7").Value = "'61.72"
$ws.Range("E7").Value = "  +13.80%  "

$ws.Range("E8").Value = "  -0.13%  "

$ws.Range("D9").Value = "'58.76"
$ws.Range("E9").Value = "  -0.44%  "

$ws.Range("D10").Value = "'0.371"
$ws.Range("E10").Value = "  +2.17%  "

$ws.Range("D11").Value = "'0.0744"
$ws.Range("E11").Value = "  +1.11%  "

$ws.Range("E12").Value = "  -1.67%  "

$ws.Range("D13").Value = "'0.894"
$ws.Range("E13").Value = "  -0.02%  "

$ws.Range("D14").Value = "'14.89"
$ws.Range("E14").Value = "  +5.19%  "

$ws.Range("D15").Value = "'2.305.23"
$ws.Range("E15").Value = "  -1.36%  "

$ws.Range("D16").Value = "'20.41"
$ws.Range("E16").Value = "  +17.57%  "

$ws.Range("D17").Value = "'5.45"
$ws.Range("E17").Value = "  +2.64%  "

$ws.Range("D18").Value = "'2.007.81"
$ws.Range("E18").Value = "  -1.37%  "

$ws.Range("D19").Value = "'36.396.45"
$ws.Range("E19").Value = "  +0.76%  "

$ws.Range("D20").Value = "'72.06"
$ws.Range("E20").Value = "  +1.22%  "

$ws.Range("D21").Value = "'0.0₃0864"
$ws.Range("E21").Value = "  +1.71%  "

$ws.Range("D22").Value = "'5.26"
$ws.Range("E22").Value = "  +2.01%  "

$ws.Range("D23").Value = "'234.52"
$ws.Range("E23").Value = "  -0.50%  "

$ws.Range("D24").Value = "'2.76"
$ws.Range("E24").Value = "  +21.21%  "

$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("D26").Value = "'2.33"
$ws.Range("E26").Value = "  -0.53%  "

$ws.Range("D27").Value = "'9.55"
$ws.Range("E27").Value = "  +4.10%  "

$ws.Range("D28").Value = "'163.88"
$ws.Range("E28").Value = "  +0.53%  "

$ws.Range("D29").Value = "'19.63"
$ws.Range("E29").Value = "  -0.89%  "

$ws.Range("E30").Value = "  -0.42%  "

$ws.Range("D31").Value = "'5.09"
$ws.Range("E31").Value = "  +3.41%  "

$ws.Range("D32").Value = "'0.110"
$ws.Range("E32").Value = "  +22.82%  "

$ws.Range("E33").Value = "  +1.78%  "

$ws.Range("D34").Value = "'4.65"
$ws.Range("E34").Value = "  +7.23%  "

$ws.Range("D35").Value = "'0.0608"
$ws.Range("E35").Value = "  +2.22%  "

$ws.Range("D36").Value = "'2.45"
$ws.Range("E36").Value = "  +11.57%  "

$ws.Range("E37").Value = "  -0.21%  "

$ws.Range("E38").Value = "  -1.12%  "

$ws.Range("D39").Value = "'5.90"
$ws.Range("E39").Value = "  +17.39%  "

$ws.Range("E40").Value = "  +15.12%  "

$ws.Range("D41").Value = "'2.80"
$ws.Range("E41").Value = "  +24.58%  "

$ws.Range("D42").Value = "'1.22"
$ws.Range("E42").Value = "  +2.10%  "

$ws.Range("D43").Value = "'2.93"
$ws.Range("E43").Value = "  +0.91%  "

$ws.Range("E44").Value = "  +3.39%  "

$ws.Range("D45").Value = "'0.0216"
$ws.Range("E45").Value = "  +0.68%  "

$ws.Range("E46").Value = "  +7.71%  "

$ws.Range("E47").Value = "  +8.76%  "

$ws.Range("D48").Value = "'94.51"
$ws.Range("E48").Value = "  +3.03%  "

$ws.Range("D49").Value = "'1.428.80"
$ws.Range("E49").Value = "  +2.26%  "

$ws.Range("D50").Value = "'2.92"
$ws.Range("E50").Value = "  -0.57%  "

$ws.Range("D51").Value = "'47.04"
$ws.Range("E51").Value = "  +3.33%  "
